$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "JAIME EMANUEL ES" rows (old rows 2 and 3) are removed.
# The remaining rows (old rows 4 and 5) shift up to become rows 2 and 3.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
